# Update the "Förändrad" (Changed) date column (C) for rows 2-16
# from 45233 (2023-11-03) to 45243 (2023-11-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value2 = 45243
    }
}
